# Updates cryptos list values (Price and Volume(1h) columns) per the
# commit "Updated cryptos list on Tue Oct 24 17:56:28 UTC 2023 with GitHub Actions".
# Rows 35/36 (Maker / LidoDAOToken) also swap position in the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume columns as plain text so values such as "3.80" or
# "1.80" are not coerced into numbers (which would drop trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "33.901.01"
$ws.Range("E2").Value = "  +9.30%  "
$ws.Range("D3").Value = "1.781.17"
$ws.Range("E3").Value = "  +5.59%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "224.96"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +4.35%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "30.61"
$ws.Range("E8").Value = "  +4.08%  "
$ws.Range("D9").Value = "46.54"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "0.277"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").Value = "0.0924"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "2.037.79"
$ws.Range("E13").Value = "  +5.63%  "
$ws.Range("D14").Value = "1.787.33"
$ws.Range("E14").Value = "  +6.12%  "
$ws.Range("D15").Value = "0.627"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "33.876.70"
$ws.Range("E16").Value = "  +9.10%  "
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").Value = "250.99"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "0.0₃0738"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "2.15"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "159.04"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "16.47"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "6.93"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").Value = "  +7.90%  "
$ws.Range("D32").Value = "0.0513"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.80"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.481.69"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").Value = "0.631"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").Value = "83.33"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +5.21%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").Value = "0.0507"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "1.937.07"
$ws.Range("E47").Value = "  +6.43%  "
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "11.97"
$ws.Range("E50").Value = "  +15.85%  "
$ws.Range("D51").Value = "50.66"
$ws.Range("E51").Value = "  -2.63%  "
